$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = "MCH239"
$ws.Range("C2").Value = "MARXIST THEORY SEMINAR GROUP 1991 CONFERENCE"
$ws.Range("D2").Value = "1991"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1C | GRAP COUNT NUMER: NONE"
